$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 463.732605
$ws.Range("H2").Value = 1391.197815
$ws.Range("I2").Value = 0.3632113435366598
$ws.Range("J2").Value = 0.3632113435366598
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 239.0839323333333
$ws.Range("N2").Value = 717.251797
$ws.Range("O2").Value = 0.4086975387666237
$ws.Range("P2").Value = 0.4086975387666237
$ws.Range("Q2").Value = 110871.0147545804
$ws.Range("R2").Value = 997839.1327912236
$ws.Range("S2").Value = 0.1484435821555515
$ws.Range("T2").Value = 0.1484435821555515

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 463.732605
$ws.Range("H3").Value = 1391.197815
$ws.Range("I3").Value = 0.3632113435366598
$ws.Range("J3").Value = 0.3632113435366598
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 117.0512696666667
$ws.Range("N3").Value = 351.153809
$ws.Range("O3").Value = 0.2000910950200451
$ws.Range("P3").Value = 0.2000910950200451
$ws.Range("Q3").Value = 54280.49020108081
$ws.Range("R3").Value = 488524.4118097274
$ws.Range("S3").Value = 0.07267535545195204
$ws.Range("T3").Value = 0.07267535545195204

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 463.732605
$ws.Range("H4").Value = 1391.197815
$ws.Range("I4").Value = 0.3632113435366598
$ws.Range("J4").Value = 0.3632113435366598
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 171.15883
$ws.Range("N4").Value = 513.47649
$ws.Range("O4").Value = 0.2925842480357353
$ws.Range("P4").Value = 0.2925842480357353
$ws.Range("Q4").Value = 79371.93010465214
$ws.Range("R4").Value = 714347.3709418693
$ws.Range("S4").Value = 0.1062699178267228
$ws.Range("T4").Value = 0.1062699178267228

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 463.732605
$ws.Range("H5").Value = 1391.197815
$ws.Range("I5").Value = 0.3632113435366598
$ws.Range("J5").Value = 0.3632113435366598
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 57.695868
$ws.Range("N5").Value = 173.087604
$ws.Range("O5").Value = 0.09862711817759588
$ws.Range("P5").Value = 0.09862711817759588
$ws.Range("Q5").Value = 26755.45516537614
$ws.Range("R5").Value = 240799.0964883853
$ws.Range("S5").Value = 0.03582248810243353
$ws.Range("T5").Value = 0.03582248810243353

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.24662533333333
$ws.Range("H6").Value = 36.739876
$ws.Range("I6").Value = 0.009591978638444229
$ws.Range("J6").Value = 0.009591978638444227
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 239.0839323333333
$ws.Range("N6").Value = 717.251797
$ws.Range("O6").Value = 0.4086975387666237
$ws.Range("P6").Value = 0.4086975387666237
$ws.Range("Q6").Value = 2927.971342506353
$ws.Range("R6").Value = 26351.74208255718
$ws.Range("S6").Value = 0.003920218061434186
$ws.Range("T6").Value = 0.003920218061434185

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.24662533333333
$ws.Range("H7").Value = 36.739876
$ws.Range("I7").Value = 0.009591978638444229
$ws.Range("J7").Value = 0.009591978638444227
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 117.0512696666667
$ws.Range("N7").Value = 351.153809
$ws.Range("O7").Value = 0.2000910950200451
$ws.Range("P7").Value = 0.2000910950200451
$ws.Range("Q7").Value = 1433.483044398632
$ws.Range("R7").Value = 12901.34739958769
$ws.Range("S7").Value = 0.001919269509175187
$ws.Range("T7").Value = 0.001919269509175187

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.24662533333333
$ws.Range("H8").Value = 36.739876
$ws.Range("I8").Value = 0.009591978638444229
$ws.Range("J8").Value = 0.009591978638444227
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 171.15883
$ws.Range("N8").Value = 513.47649
$ws.Range("O8").Value = 0.2925842480357353
$ws.Range("P8").Value = 0.2925842480357353
$ws.Range("Q8").Value = 2096.118063501694
$ws.Range("R8").Value = 18865.06257151524
$ws.Range("S8").Value = 0.002806461857104041
$ws.Range("T8").Value = 0.00280646185710404

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.24662533333333
$ws.Range("H9").Value = 36.739876
$ws.Range("I9").Value = 0.009591978638444229
$ws.Range("J9").Value = 0.009591978638444227
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 57.695868
$ws.Range("N9").Value = 173.087604
$ws.Range("O9").Value = 0.09862711817759588
$ws.Range("P9").Value = 0.09862711817759588
$ws.Range("Q9").Value = 706.5796786774561
$ws.Range("R9").Value = 6359.217108097105
$ws.Range("S9").Value = 0.0009460292107308142
$ws.Range("T9").Value = 0.000946029210730814

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 689.7685036666667
$ws.Range("H10").Value = 2069.305511
$ws.Range("I10").Value = 0.5402504422695089
$ws.Range("J10").Value = 0.5402504422695089
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 239.0839323333333
$ws.Range("N10").Value = 717.251797
$ws.Range("O10").Value = 0.4086975387666237
$ws.Range("P10").Value = 0.4086975387666237
$ws.Range("Q10").Value = 164912.5662563059
$ws.Range("R10").Value = 1484213.096306753
$ws.Range("S10").Value = 0.2207990260731282
$ws.Range("T10").Value = 0.2207990260731282

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 689.7685036666667
$ws.Range("H11").Value = 2069.305511
$ws.Range("I11").Value = 0.5402504422695089
$ws.Range("J11").Value = 0.5402504422695089
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 117.0512696666667
$ws.Range("N11").Value = 351.153809
$ws.Range("O11").Value = 0.2000910950200451
$ws.Range("P11").Value = 0.2000910950200451
$ws.Range("Q11").Value = 80738.27913026015
$ws.Range("R11").Value = 726644.5121723415
$ws.Range("S11").Value = 0.1080993025787697
$ws.Range("T11").Value = 0.1080993025787697

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 689.7685036666667
$ws.Range("H12").Value = 2069.305511
$ws.Range("I12").Value = 0.5402504422695089
$ws.Range("J12").Value = 0.5402504422695089
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 171.15883
$ws.Range("N12").Value = 513.47649
$ws.Range("O12").Value = 0.2925842480357353
$ws.Range("P12").Value = 0.2925842480357353
$ws.Range("Q12").Value = 118059.9700584374
$ws.Range("R12").Value = 1062539.730525936
$ws.Range("S12").Value = 0.1580687694023977
$ws.Range("T12").Value = 0.1580687694023977

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 689.7685036666667
$ws.Range("H13").Value = 2069.305511
$ws.Range("I13").Value = 0.5402504422695089
$ws.Range("J13").Value = 0.5402504422695089
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 57.695868
$ws.Range("N13").Value = 173.087604
$ws.Range("O13").Value = 0.09862711817759588
$ws.Range("P13").Value = 0.09862711817759588
$ws.Range("Q13").Value = 39796.79253810951
$ws.Range("R13").Value = 358171.1328429856
$ws.Range("S13").Value = 0.05328334421521329
$ws.Range("T13").Value = 0.05328334421521329

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 111.00921
$ws.Range("H14").Value = 333.02763
$ws.Range("I14").Value = 0.08694623555538696
$ws.Range("J14").Value = 0.08694623555538696
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 239.0839323333333
$ws.Range("N14").Value = 717.251797
$ws.Range("O14").Value = 0.4086975387666237
$ws.Range("P14").Value = 0.4086975387666237
$ws.Range("Q14").Value = 26540.51845201679
$ws.Range("R14").Value = 238864.6660681511
$ws.Range("S14").Value = 0.03553471247650975
$ws.Range("T14").Value = 0.03553471247650975

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 111.00921
$ws.Range("H15").Value = 333.02763
$ws.Range("I15").Value = 0.08694623555538696
$ws.Range("J15").Value = 0.08694623555538696
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 117.0512696666667
$ws.Range("N15").Value = 351.153809
$ws.Range("O15").Value = 0.2000910950200451
$ws.Range("P15").Value = 0.2000910950200451
$ws.Range("Q15").Value = 12993.76897519363
$ws.Range("R15").Value = 116943.9207767427
$ws.Range("S15").Value = 0.01739716748014815
$ws.Range("T15").Value = 0.01739716748014815

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 111.00921
$ws.Range("H16").Value = 333.02763
$ws.Range("I16").Value = 0.08694623555538696
$ws.Range("J16").Value = 0.08694623555538696
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 171.15883
$ws.Range("N16").Value = 513.47649
$ws.Range("O16").Value = 0.2925842480357353
$ws.Range("P16").Value = 0.2925842480357353
$ws.Range("Q16").Value = 19000.2065028243
$ws.Range("R16").Value = 171001.8585254187
$ws.Range("S16").Value = 0.0254390989495108
$ws.Range("T16").Value = 0.0254390989495108

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 111.00921
$ws.Range("H17").Value = 333.02763
$ws.Range("I17").Value = 0.08694623555538696
$ws.Range("J17").Value = 0.08694623555538696
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 57.695868
$ws.Range("N17").Value = 173.087604
$ws.Range("O17").Value = 0.09862711817759588
$ws.Range("P17").Value = 0.09862711817759588
$ws.Range("Q17").Value = 6404.772726944279
$ws.Range("R17").Value = 57642.95454249852
$ws.Range("S17").Value = 0.008575256649218238
$ws.Range("T17").Value = 0.008575256649218238
